$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.6606524410359556, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 3.56341032713086),
    @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694),
    @(0.2917716402565462, 1.655778082260271, 3.537761648806719, 1133.036916526867, 1138.522227898191),
    @(3.286832544864788, 1.655778082260271, 3.537761648806719, 10.19245300693656, 18.67282528286833),
    @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 10.19245300693656, 15.88780690183548),
    @(0.2917716402565462, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 2.591208233317391),
    @(0.2917716402565462, 0.306821227259698, 0.1494219747398047, 0.4942365360607697, 1.242251378316819),
    @(0.2917716402565462, 0.306821227259698, 3.537761648806719, 0.4942365360607697, 4.630591052383734),
    @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 5.586269137925634),
    @(3.286832544864788, 1.655778082260271, 3.537761648806719, 10.19245300693656, 18.67282528286833),
    @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694),
    @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 5.586269137925634),
    @(1.455362044514542, 1.655778082260271, 22.3905356188092, 0.4942365360607697, 25.99591228164478),
    @(0.1190320826869504, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 2.418468675747795),
    @(1.455362044514542, 1.655778082260271, 3.537761648806719, 10.19245300693656, 16.84135478251809),
    @(1.455362044514542, 0.04071648406533734, 0.7527432677738641, 0.4942365360607697, 2.743058332414513),
    @(0.6606524410359556, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 3.56341032713086),
    @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694),
    @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694),
    @(3.286832544864788, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 5.586269137925634),
    @(3.286832544864788, 1.655778082260271, 3.537761648806719, 10.19245300693656, 18.67282528286833),
    @(0.6606524410359556, 0.04071648406533734, 0.7527432677738641, 0.4942365360607697, 1.948348728935927),
    @(1.455362044514542, 1.655778082260271, 0.1494219747398047, 0.4942365360607697, 3.754798637575387),
    @(0.6606524410359556, 0.306821227259698, 0.1494219747398047, 0.4942365360607697, 1.611132179096228),
    @(3.286832544864788, 1.655778082260271, 0.7527432677738641, 0.4942365360607697, 6.189590430959694),
    @(1.455362044514542, 1.655778082260271, 0.7527432677738641, 10.19245300693656, 14.05633640148523),
    @(3.286832544864788, 1.655778082260271, 3.537761648806719, 0.4942365360607697, 8.974608811992548)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 2).Value = $data[$i][0]
    $ws.Cells.Item($row, 3).Value = $data[$i][1]
    $ws.Cells.Item($row, 4).Value = $data[$i][2]
    $ws.Cells.Item($row, 5).Value = $data[$i][3]
    $ws.Cells.Item($row, 7).Value = $data[$i][4]
}

$wb.Save()